$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newly discovered facility names to the bottom of the list
# (column A), continuing the existing single-column table.
$ws.Range("A67").Value = "LUDLOW COMMUNITY CENTER"
$ws.Range("A68").Value = "Bluford Charter School"
$ws.Range("A69").Value = "Duckery School"
$ws.Range("A70").Value = "Building 21 @ Kinsey School"
$ws.Range("A71").Value = "Pennell School Kindergarten Mobile bldg"

# Update the visible selection to match where the new rows were added.
$ws.Range("A69:A71").Select()
